$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: refreshed timestamp banner
$ws.Range("A1").Value = 'Last updated: 2025-07-15 07:58:12'

# Data rows 3-29: refreshed PO status query results
$data = @(
    @{ Row=3; A='4516260169'; B=4; C=0; D=0; E=0; F='Completed PO'; G=0; H=4; I=0 },
    @{ Row=4; A='47188646'; B=2; C=0; D=2; E=0; F='Completed PO'; G=0; H=0; I=0 },
    @{ Row=5; A='3T/PO252272'; B=1; C=0; D=0; E=0; F='Completed PO'; G=1; H=0; I=0 },
    @{ Row=6; A='4010016033'; B=1; C=0; D=0; E=0; F='Completed PO'; G=1; H=0; I=0 },
    @{ Row=7; A='4020007186'; B=11; C=0; D=0; E=0; F='Completed PO'; G=11; H=0; I=0 },
    @{ Row=8; A='4516260169'; B=3; C=0; D=0; E=0; F='Completed PO'; G=3; H=0; I=0 },
    @{ Row=9; A='4516351202_AIZU'; B=15; C=2; D=1; E=0; F=$null; G=12; H=0; I=-1 },
    @{ Row=10; A='4516351202_ARD'; B=25; C=1; D=8; E=0; F=$null; G=16; H=0; I=-1 },
    @{ Row=11; A='4516351202_ASEWH'; B=14; C=0; D=0; E=0; F='Completed PO'; G=14; H=0; I=0 },
    @{ Row=12; A='4516351202_ATK'; B=19; C=0; D=0; E=0; F='Completed PO'; G=19; H=0; I=0 },
    @{ Row=13; A='4516351202_DMOS5'; B=2; C=0; D=0; E=0; F='Completed PO'; G=2; H=0; I=0 },
    @{ Row=14; A='4516351202_HNT'; B=12; C=0; D=0; E=0; F='Completed PO'; G=12; H=0; I=0 },
    @{ Row=15; A='4516351202_LFAB'; B=1; C=0; D=0; E=1; F='Completed PO'; G=0; H=0; I=0 },
    @{ Row=16; A='4516351202_MIHO'; B=31; C=0; D=0; E=0; F='Completed PO'; G=31; H=0; I=0 },
    @{ Row=17; A='4516351202_SCT'; B=29; C=0; D=0; E=0; F='Completed PO'; G=29; H=0; I=0 },
    @{ Row=18; A='4516351202_TICL-FT'; B=31; C=3; D=7; E=0; F=$null; G=21; H=0; I=-1 },
    @{ Row=19; A='4516351202_TICL-PR'; B=6; C=0; D=0; E=1; F='Completed PO'; G=5; H=0; I=0 },
    @{ Row=20; A='4516351202_TII'; B=8; C=0; D=0; E=0; F='Completed PO'; G=8; H=0; I=0 },
    @{ Row=21; A='4516351202_TIPI'; B=51; C=5; D=8; E=10; F=$null; G=28; H=0; I=-1 },
    @{ Row=22; A='4516351202_TITL'; B=56; C=0; D=3; E=0; F=$null; G=53; H=0; I=0 },
    @{ Row=23; A='4516351202_UTL'; B=1; C=0; D=0; E=0; F='Completed PO'; G=1; H=0; I=0 },
    @{ Row=24; A='47188646'; B=2; C=0; D=0; E=0; F='Completed PO'; G=2; H=0; I=0 },
    @{ Row=25; A='47225672'; B=1; C=0; D=0; E=1; F='Completed PO'; G=0; H=0; I=0 },
    @{ Row=26; A='4910171717'; B=10; C=0; D=2; E=6; F='Completed PO'; G=2; H=0; I=0 },
    @{ Row=27; A='Global PO 4516351202'; B=19; C=0; D=7; E=9; F=$null; G=3; H=0; I=0 },
    @{ Row=28; A='PP25-1057NI'; B=3; C=0; D=0; E=0; F='Completed PO'; G=3; H=0; I=0 },
    @{ Row=29; A='TBA'; B=24; C=0; D=14; E=1; F=$null; G=9; H=0; I=0 }
)

foreach ($r in $data) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    if ($null -eq $r.F) {
        $ws.Cells.Item($r.Row, 6).Value = ""
    } else {
        $ws.Cells.Item($r.Row, 6).Value = $r.F
    }
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
    $ws.Cells.Item($r.Row, 9).Value = $r.I
}

Write-Output "PO Status updated"
